$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.630.42'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +8.59%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.496.64'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +11.89%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.22%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '188.97'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +13.09%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '549.67'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +8.32%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.488.03'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +11.77%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.06%  '

# Row 9
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.01%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.632'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.70%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.150'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +19.26%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.17'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +8.23%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +9.57%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.36'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +8.80%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.071.07'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +10.91%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.512.43'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +11.34%  '

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +8.22%  '

# Row 18
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.777.88'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +8.79%  '

# Row 19
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.20'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +9.55%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.76'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +11.52%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.990'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +6.65%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '414.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +17.02%  '

# Row 23
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '85.18'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.77%  '

# Row 24
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.91'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +8.67%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.21'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +12.57%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.16'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.90%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.92'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +16.30%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.14'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.96%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.85'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +9.92%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.81'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +12.16%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.20'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +10.26%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '654.98'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.78%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.69'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +8.26%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.68'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.65%  '

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +10.32%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.49'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.37%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.68'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +10.28%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0809'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +20.77%  '

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.19%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.391'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.68%  '

# Row 41
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +13.94%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.36'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +22.68%  '

# Row 43
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.07%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.992.34'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +7.38%  '

# Row 45
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.63'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +9.69%  '

# Row 46
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.91'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +17.94%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.34'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +17.45%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0416'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +11.62%  '

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.46%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.93'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +22.35%  '

# Row 51
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +9.22%  '
